# aug 2 update 4 appendix a
# Add two new Skeena commercial fishery opening rows (FN0762, FN0763) to the
# "Skeena" sheet, then leave a blank (but date-formatted) cell below them,
# matching row 11's quota-style number formatting for the Aug 2-5 / Sockeye row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Skeena")

# Row 22: FN0762 / Commercial / Aug 2-3 / Sockeye / Gillnet / Area 4 / 2 days / Pink retention
$ws.Range("A22").Value = "FN0762"
$ws.Range("B22").Value = "Commercial"
$ws.Range("C22").Value = "Aug 2-3"
$ws.Range("D22").Value = "Sockeye"
$ws.Range("E22").Value = "Gillnet"
$ws.Range("F22").Value = "Area 4"
$ws.Range("G22").Value = 2
$ws.Range("I22").Value = "Pink retention"

# Row 23: FN0763 / Commercial / Aug 2-5 / Sockeye / Seine / Area 4 / 4 days / quota 46656 / Pink retention
# C23 and D23 reuse the quote-prefixed date-ish text style seen on row 11 (C11/D11),
# so force-text (leading apostrophe) before applying the "d-mmm" number format.
$ws.Range("A23").Value = "FN0763"
$ws.Range("B23").Value = "Commercial"
$ws.Range("C23").Value = "'Aug 2-5"
$ws.Range("C23").NumberFormat = "d-mmm"
$ws.Range("D23").Value = "'Sockeye"
$ws.Range("D23").NumberFormat = "d-mmm"
$ws.Range("E23").Value = "Seine"
$ws.Range("F23").Value = "Area 4"
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = 46656
$ws.Range("I23").Value = "Pink retention"

# Row 24: a single formatted-but-empty cell at C24 (same numeric date format, no quote prefix)
$ws.Range("C24").Value = ""
$ws.Range("C24").NumberFormat = "d-mmm"

# Update the active selection to C12, matching the saved sheet view state.
$ws.Range("C12").Select()
